# Insert a new row of localization data ("launch" / "LAUNCH") above the
# existing "newton_first_law_dlg_1" row, pushing all following rows down
# by one. This mirrors Excel's native Insert-Row behaviour, which also
# renumbers the sparse block of styled-but-empty rows further down the
# sheet (the act-2 placeholder geometry rows).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new row 31 (shifts row 31.. down to row 32..).
$ws.Rows(31).Insert()

# Populate the newly inserted row with the new localization key/value pair.
$ws.Cells.Item(31, 1).Value = "launch"
$ws.Cells.Item(31, 2).Value = "LAUNCH"

# Match the author's recorded selection/viewport after the edit.
$ws.Range("B31").Select()
